$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# L2 ("Date" column) currently holds the plain number 2018; change it to a
# real Excel date serial value and format it as a proper date (DD/MM/YY).
$ws.Range("L2").Value = 43185
$ws.Range("L2").NumberFormat = "DD/MM/YY"

# Incidental change captured by the diff: the active selection ends up on L3
# (the cell below the one that was just edited).
$ws.Range("L3").Select()
